$d = $word.ActiveDocument

# Replace the first occurrence of $find (plain text) with $replace, without
# letting the new text inherit formatting from a neighbouring run (which
# happens if we let Find.Execute's own Replace do the splice when the match
# starts exactly at a run/hyperlink boundary). We locate the match, insert
# the new text immediately after the match (a position that is always safely
# inside -- or at the very end of -- the original run), and only then delete
# the original characters.
function Replace-Text($find, $replace) {
    $r = $d.Content
    $ok = $r.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        Write-Output "WARNING: text not found: $find"
        return
    }
    $s = $r.Start
    $e = $r.End
    $ins = $d.Range($e, $e)
    $ins.InsertBefore($replace)
    $old = $d.Range($s, $e)
    $old.Delete()
}

# 1. Intro paragraph: remove duplicated "para"
Replace-Text "Este guia contém as orientações básicas para sobre o processo de " "Este guia contém as orientações básicas sobre o processo de "

# 2. "O processo ..." paragraph: rewording + "ambiente" + italic Kubernetes
Replace-Text " para levar a aplicação desde o seu código fonte inicial até sua execução final no " " para levar a aplicação desde o seu código fonte inicial até a sua execução final no ambiente "

$r = $d.Content
$r.Find.Execute("Kubernetes requer")
$kw = $d.Range($r.Start, $r.Start + 10)
$kw.Font.Italic = $true

Replace-Text "uma série de etapas a serem compridas." "um conjunto de etapas que precisam ser cumpridas."

# 3. "Mapeando..." paragraph reworded
Replace-Text "Mapeando os principais pontos desse processo podemos dividi-lo da seguinte forma:" "Fazendo um mapeando dos principais pontos desse processo podemos dividi-lo da seguinte maneira:"

# 4. "Segue abaixo..." paragraph reworded
Replace-Text "Segue abaixo o detalhamento de cada umas das etapas." "Tomando como base as 5 etapas acima, segue abaixo um detalhamento de cada uma delas."

# 5. "Para mais detalhes..." paragraph reworded
Replace-Text "segue abaixo algumas referências sobre o assunto:" "sobre o assunto, seguem abaixo algumas referências:"

# 6. "... de uma aplicação." -> "... da aplicação." (only first occurrence)
Replace-Text " de uma aplicação" " da aplicação"

# 7. "Para teste e desenvolvimento local..." paragraph reworded
Replace-Text "Para teste e desenvolvimento local é necessário ter instalado a ferramenta do Docker. " "Para realizar testes durante o desenvolvimento local do arquivo, é necessário ter instalado a ferramenta do Docker. "

# 8. "... para obter mais detalhes." -> add "sobre a instalação"
Replace-Text " para obter mais detalhes." " para obter mais detalhes sobre a instalação."

# 9. "Dockerfile, basicamente são necessários apenas 2 comandos" reworded
Replace-Text ", basicamente são necessários apenas 2 comandos" " no diretório raiz da aplicação, basicamente, serão necessários apenas 2 comandos"

# 10. "A referência completa de comandos para o Docker CLI" reworded
Replace-Text "A referência completa de comandos para o Docker CLI" "A referência completa para os comandos do Docker CLI"

# 11. Move the _GoBack bookmark from its old standalone paragraph to right
# after the "." that ends the "A referência completa..." sentence (before
# the following line break run). Search without the trailing "." (which
# sits in its own following run) and skip the two zero-width marks (the
# hyperlink end + the "." run) to land right after the period.
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

$r = $d.Content
$r.Find.Execute("pode ser consultada aqui")
$bmPos = $r.End + 2
$d.Bookmarks.Add("_GoBack", $d.Range($bmPos, $bmPos)) | Out-Null
